# Remove Sheet2 entirely (workbook entry, worksheet part, and its shared-string usages go with it)
$wb = $excel.ActiveWorkbook
$wb.Worksheets("Sheet2").Delete()

# Update Sheet1 values
$ws = $wb.Worksheets("Sheet1")
$ws.Range("C1").Value = 24
$ws.Range("C2").Value = 25
$ws.Range("B3").Value = 35
$ws.Range("C3").Value = 21
$ws.Range("B4").Value = 44
$ws.Range("C4").Value = 24

# Update the active selection on Sheet1
$ws.Range("D9").Select()
